# The brand-name column (B) is replaced with numeric figures.
# Column A (tea types) stays the same text, only column B changes
# from shared-string brand names to plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, New B value
$values = @(
    , @(2, 10)
    , @(3, 30)
    , @(4, 20)
    , @(5, 23)
    , @(6, 54)
    , @(7, 23)
    , @(8, 14)
    , @(9, 16)
    , @(10, 18)
    , @(11, 45)
    , @(12, 34)
    , @(13, 45)
    , @(14, 36)
)

foreach ($pair in $values) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 2).Value = $val
}

# Restore the current selection/viewport state saved with the workbook.
$ws.Range("B14").Select()
